$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(13.38179432741989,9.821207075173829,15.11605234562326,16.54754509988564,49.82800970855569,19.46201495482082,27.94322330669623,9.461007956457282,14.2725942627142,21.27912208954055)
    3 = @(13.16826753186891,9.666647065647995,15.05985833177668,16.49197439685554,49.75429207266262,19.49564544420134,27.9946642613009,9.470674058074998,14.13226323546643,21.34111704883852)
    4 = @(13.03891796896394,9.573317528262857,15.02874844705138,16.46159628712114,49.72189924423194,19.51961667842892,28.0314211882182,9.478064807865408,14.04889636946693,21.38108631688144)
    5 = @(12.98672042628981,9.535731069025626,15.01693301909318,16.45016628363581,49.71193900943327,19.53021924043901,28.0476976660889,9.481442648853584,14.01566307864856,21.39785388777137)
    6 = @(12.97808616108288,9.529518287474982,15.01502338113051,16.44832589589446,49.71048086848104,19.53203012584893,28.05047865657634,9.482025643817838,14.01019039764888,21.40066713560349)
    7 = @(13.03821184133579,9.572808750987333,15.02858559861631,16.46143828426771,49.72175179523244,19.51975629288252,28.03163544795392,9.478108880546163,14.04844513342764,21.3813105063867)
    8 = @(13.30784879736589,9.767619540485031,15.09597792318332,16.52761282910087,49.79992382841471,19.4729205529598,27.95988507780245,9.464038709866442,14.22364969478401,21.3001033354688)
    9 = @(13.84738921259609,10.15984522275789,15.254626493454,16.68666548146455,50.05504464462312,19.40748675069674,27.86033889999142,9.447996026827973,14.58776768072216,21.15591973487927)
    10 = @(14.24617627047738,10.45121263485618,15.3866901468615,16.82075529777097,50.30398321372026,19.37558021312326,27.81244478808573,9.443245051793406,14.86552480499259,21.05910717686086)
    11 = @(14.42725162762469,10.5838315514214,15.44997857558337,16.88534542754508,50.43040202523132,19.36458746058991,27.79616568421258,9.442608846861743,14.99360776174957,21.01703134016554)
    12 = @(14.49570292811852,10.63401107575497,15.47439148782505,16.91030614704051,50.48014767080872,19.36093188446038,27.79079509656031,9.442586844540852,15.04231596496779,21.00137989335181)
    13 = @(14.48096704816838,10.62320660787971,15.46911409780287,16.90490833308327,50.46935109315061,19.36169661306526,27.79191641017301,9.442581853884201,15.0318173011259,21.00473819300112)
    14 = @(14.4328859130556,10.58796094672831,15.45197817115405,16.88738899055044,50.43445721067583,19.36427654668875,27.79570792085074,9.442602651849597,14.99761113162658,21.01573804361726)
    15 = @(14.403417433249,10.56636521652736,15.44153967193721,16.87672279097873,50.41332702294169,19.36592289407941,27.79813377840045,9.442643887479459,14.97668442052581,21.02251244093173)
    16 = @(14.23432980123463,10.44254276856432,15.3826174755033,16.81660522576578,50.29598483926356,19.37636954637834,27.81361967187065,9.443317286106907,14.85718543231447,21.0618964007625)
    17 = @(14.13046373774959,10.36656389039275,15.34728292716105,16.78063487644899,50.22735971686348,19.3836807845707,27.82453197002537,9.444120747487489,14.78428927106039,21.08655983968339)
    18 = @(14.07069523201189,10.32287281147158,15.32726270567236,16.76028471594083,50.18912993839595,19.38821745985064,27.83132679417266,9.444726449410263,14.7425267119783,21.10093056651883)
    19 = @(14.05045609141281,10.30808305811154,15.32053673220333,16.75345316351086,50.17639982729614,19.38981040145113,27.83371637158019,9.444956198611445,14.72841633683518,21.10582803984554)
    20 = @(14.14152381719858,10.37465131611806,15.35101306779477,16.78442899650336,50.23453662269961,19.38286818172868,27.82331667444976,9.444020360972292,14.79203238535906,21.08391523642799)
    21 = @(14.44701225419333,10.59831494608864,15.45699940073597,16.8925213476872,50.44465572307969,19.36350498898129,27.79457270084725,9.442590605192262,15.00765306317077,21.01249947992943)
    22 = @(14.64594830217395,10.74423428596093,15.52886534525599,16.96608364264464,50.59288962093104,19.35380625269463,27.78041531328347,9.442931939714697,15.1497544789511,20.96746712635257)
    23 = @(14.53986049298386,10.66639438381345,15.4902765746073,16.9265601639523,50.51278391990378,19.35871195647251,27.78754732103713,9.44263318004432,15.07381831949305,20.99135174601933)
    24 = @(14.13652371971682,10.37099501757777,15.34932575406961,16.78271264684628,50.23128813003122,19.38323452109313,27.82386448611686,9.444065297820302,14.78853126401999,21.08511026527232)
    25 = @(13.70071146823354,10.0529595603479,15.20893445112801,16.64056215885693,49.97517393624462,19.42235340639459,27.8828460262331,9.451099799937733,14.48729756114536,21.19331925835288)
}

$cols = @(2,3,4,5,7,8,9,10,11,14)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}
